$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins, Losses, Ties in AD1, AE1, AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill data rows 2-51 with Wins=85, Losses=76, Ties=0
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 85  # AD
    $ws.Cells.Item($r, 31).Value = 76  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
